# Rename the sheet "C_29" to "C_39"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "C_39"

# Restore the active cell selection to B2 (top-left of the used range)
$ws.Range("B2").Select()
